# Loan product workbook cleanup:
#  - shorten the "description" field (B3) to the short product code "7A"
#  - normalise the currency label/value row (A6/B6): lower-case key,
#    drop the trailing space from the currency name
#  - drop the leftover formatting-only cells and blank trailing row
#  - tidy up the sheet view (selection / zoom) on the input sheet

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# --- content edits -------------------------------------------------
$ws1.Cells.Item(3, 2).Value = "7A"

$ws1.Cells.Item(6, 1).Value = "currency"
$ws1.Cells.Item(6, 2).Value = "US Dollar"

# --- drop stray formatting-only cells that no longer carry data ----
$ws1.Cells.Item(5, 3).Clear()
$ws1.Cells.Item(6, 3).Clear()
$ws1.Cells.Item(12, 3).Clear()
$ws1.Cells.Item(26, 4).Clear()
$ws1.Cells.Item(26, 5).Clear()

# --- drop the empty trailing row ------------------------------------
$ws1.Rows.Item(42).Delete()

# --- sheet view tidy-up ----------------------------------------------
[void]$ws1.Range("A6:B6").Select()
$excel.ActiveWindow.Zoom = 100

[void]$ws2.Cells.Item(1, 1).Select()
